# "removed literally a single dash"
#
# Slide 7 ("Which Metrics Did Homies Use?"), the metrics list text box,
# has a line that reads "Pollution \u2013" (an en dash trailing the word).
# The author simply deleted the " \u2013" so the line reads "Pollution".

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(7)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange

$dash   = [char]0x2013
$needle = "Pollution " + $dash

$hit = $tr.Find($needle, 0)
if ($hit -ne $null) {
    $hit.Text = "Pollution"
}
